$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update existing row
$ws.Range("B2").Value = "Flu B PCR"
$ws.Range("C2").Value = 4008926
$ws.Range("D2").Value = 2023049001139
$ws.Range("E2").Value = "POSITIVE"

# New rows 3-9
$data = @(
    @(1, "Flu A PCR", 4001509, 2023046001696, "POSITIVE", "Categorical"),
    @(2, "RSV PCR", 3995090, 2023044002184, "POSITIVE", "Categorical"),
    @(3, "SARSCoV2 PCR", 3994204, 2023044001487, "POSITIVE", "Categorical"),
    @(4, "SARS-CoV-2 (In-house)", 3994149, 2023043001779, "Positive", "Categorical"),
    @(5, "SARSCoV2 PCR", 3994135, 2023043001778, "Not Detected", "Categorical"),
    @(6, "SARS-CoV-2 (In-house)", 3993717, 2023044000622, "Negative", "Categorical"),
    @(7, "SARSCoV2 PCR", 3993716, 2023044000836, "NEGATIVE", "Categorical")
)

$row = 3
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $row++
}

$ws.Range("A2").Copy()
$ws.Range("A3:A9").PasteSpecial(-4122)

